# Update docx golden tests for reference doc changes.
#
# Mirrors the golden diff for test/docx/golden/tables.docx:
#   * Title / TitleChar run properties drop their explicit character
#     spacing (-10) and minimum-kerning (28) overrides.
#   * Author / Date paragraph styles now inherit from Title (instead of
#     duplicating its centering via an explicit jc) and pick up an
#     explicit 12pt run size (sz/szCs 24 half-points) in their own rPr.

$d = $word.ActiveDocument

# --- Title -----------------------------------------------------------
$title = $d.Styles("Title")
$title.Font.Spacing = 0
$title.Font.Kerning = 0

# --- TitleChar (linked character style for Title) --------------------
$titleChar = $d.Styles("TitleChar")
$titleChar.Font.Spacing = 0
$titleChar.Font.Kerning = 0

# --- Author ------------------------------------------------------------
$author = $d.Styles("Author")
$author.BaseStyle = "Title"
$author.Font.Size = 12
$author.Font.SizeBi = 12

# --- Date ----------------------------------------------------------------
$date = $d.Styles("Date")
$date.BaseStyle = "Title"
$date.Font.Size = 12
$date.Font.SizeBi = 12
